# Insert two new rows (Ciruela / Angeleno, Primera & Segunda, Región Metropolitana,
# fecha 2023-03-30) before what is currently row 102, shifting the existing rows
# 102:118 down to 104:120.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("102:103").Insert()

# New row 102
$ws.Cells.Item(102,1).Value = 2
$ws.Cells.Item(102,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(102,3).Value = "Coquimbo"
$ws.Cells.Item(102,4).Value = 45015
$ws.Cells.Item(102,5).Value = 4
$ws.Cells.Item(102,6).Value = "Fruta"
$ws.Cells.Item(102,7).Value = 100103
$ws.Cells.Item(102,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(102,9).Value = 100103002
$ws.Cells.Item(102,10).Value = "Ciruela"
$ws.Cells.Item(102,11).Value = "Angeleno"
$ws.Cells.Item(102,12).Value = "Primera"
$ws.Cells.Item(102,13).Value = 20
$ws.Cells.Item(102,14).Value = 220000
$ws.Cells.Item(102,15).Value = 230000
$ws.Cells.Item(102,16).Value = 225000
$ws.Cells.Item(102,17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(102,18).Value = "Región Metropolitana"
$ws.Cells.Item(102,19).Value = 500
$ws.Cells.Item(102,20).Value = 450

# New row 103
$ws.Cells.Item(103,1).Value = 2
$ws.Cells.Item(103,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(103,3).Value = "Coquimbo"
$ws.Cells.Item(103,4).Value = 45015
$ws.Cells.Item(103,5).Value = 4
$ws.Cells.Item(103,6).Value = "Fruta"
$ws.Cells.Item(103,7).Value = 100103
$ws.Cells.Item(103,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(103,9).Value = 100103002
$ws.Cells.Item(103,10).Value = "Ciruela"
$ws.Cells.Item(103,11).Value = "Angeleno"
$ws.Cells.Item(103,12).Value = "Segunda"
$ws.Cells.Item(103,13).Value = 14
$ws.Cells.Item(103,14).Value = 180000
$ws.Cells.Item(103,15).Value = 190000
$ws.Cells.Item(103,16).Value = 185000
$ws.Cells.Item(103,17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(103,18).Value = "Región Metropolitana"
$ws.Cells.Item(103,19).Value = 411
$ws.Cells.Item(103,20).Value = 450

# Ensure the date cells keep the original date-number-format style (style index 2)
# used by every "Fecha" cell in the sheet.
$dateFormat = $ws.Range("D104").NumberFormat
$ws.Range("D102").NumberFormat = $dateFormat
$ws.Range("D103").NumberFormat = $dateFormat
